$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to ValidLogin
$ws.Name = "ValidLogin"

# Populate test data for the ValidLogin test case.
# Order of assignment matches the shared-string table ordering in the target file.
$ws.Range("A2").Value = "admin"
$ws.Range("B1").Value = "Password"
$ws.Range("A1").Value = "Username"
$ws.Range("B2").Value = "manager"

# Update selection to B2 as in the target
$ws.Range("B2").Select()
